$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "<delete>"
$ws.Range("C2").Value = 50

# Row 3
$ws.Range("B3").Value = "<each>"
$ws.Range("C3").Value = 50

# Row 4
$ws.Range("B4").Value = "<she>"
$ws.Range("C4").Value = 52

# Row 5
$ws.Range("B5").Value = "<ou>"

# Row 6
$ws.Range("B6").Value = "<it>"
$ws.Range("C6").Value = 53

# Row 7
$ws.Range("B7").Value = "<which>"
$ws.Range("C7").Value = 53

# Row 8
$ws.Range("C8").Value = 40
